# "Three objects per class, also fixed file names"
# Update the Project Software checklist grade for the "Objects" row
# (row 25) from 1 to 2, i.e. "More than two objects per class", and
# give full credit (1) on the "Regular Expression" row (row 31) which
# was previously left blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Objects: two-or-less -> more-than-two objects per class
$ws.Range("D25").Value = 2

# Regular Expression: award the point that was previously missing
$ws.Range("D31").Value = 1

# D33 holds =SUM(D3:D32) and will recalculate automatically (27 -> 29)

# Reflect the view state recorded in the saved file: the window had been
# scrolled down so row 22 is the first visible row, with D14 selected.
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D14").Select()
